$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the header row (A1 already has "numero_facture") with the new
#     columns B..I, reusing A1's formatting (bold / centered / bordered) ---
$ws.Range("A1").Copy()
$ws.Range("B1:I1").PasteSpecial(-4122)

$ws.Range("B1").Value = "code_client"
$ws.Range("C1").Value = "nom_client"
$ws.Range("D1").Value = "date"
$ws.Range("E1").Value = "total_ht"
$ws.Range("F1").Value = "remise_fcfa"
$ws.Range("G1").Value = "taux_remise"
$ws.Range("H1").Value = "tva"
$ws.Range("I1").Value = "total_ttc"

# --- Append the new invoice record on row 17 ---
$ws.Range("A17").Value = "F0016"
$ws.Range("B17").Value = "C00003"
$ws.Range("C17").Value = "Florent b"
$ws.Range("D17").Value = "20/07/2025"
$ws.Range("E17").Value = 2380000
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 428400
$ws.Range("I17").Value = 2808400
